$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '46.501.64'
$ws.Range("E2").Value = '  +1.37%  '

$ws.Range("D3").Value = '2.543.17'
$ws.Range("E3").Value = '  +8.12%  '

$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.25%  '

$ws.Range("D5").Value = '''307.21'
$ws.Range("E5").Value = '  +2.09%  '

$ws.Range("D6").Value = '''104.04'
$ws.Range("E6").Value = '  +5.46%  '

$ws.Range("D7").Value = '''0.611'
$ws.Range("E7").Value = '  +6.90%  '

$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").Value = '''0.574'
$ws.Range("E9").Value = '  +11.32%  '

$ws.Range("D10").Value = '''39.32'
$ws.Range("E10").Value = '  +12.51%  '

$ws.Range("D11").Value = '''0.0833'
$ws.Range("E11").Value = '  +4.31%  '

$ws.Range("D12").Value = '''7.99'
$ws.Range("E12").Value = '  +11.52%  '

$ws.Range("D13").Value = '2.925.24'
$ws.Range("E13").Value = '  +7.90%  '

$ws.Range("E14").Value = '  +2.61%  '

$ws.Range("D15").Value = '2.564.00'
$ws.Range("E15").Value = '  +8.57%  '

$ws.Range("E16").Value = '  +10.11%  '

$ws.Range("D17").Value = '''15.00'
$ws.Range("E17").Value = '  +8.80%  '

$ws.Range("D18").Value = '46.563.62'
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("D19").Value = '''13.96'
$ws.Range("E19").Value = '  +10.33%  '

$ws.Range("E20").Value = '  +2.89%  '

$ws.Range("D21").Value = '''6.65'
$ws.Range("E21").Value = '  +10.83%  '

$ws.Range("D22").Value = '''70.19'
$ws.Range("E22").Value = '  +5.87%  '

$ws.Range("D23").Value = '''255.18'
$ws.Range("E23").Value = '  +4.19%  '

$ws.Range("D24").Value = '''2.99'
$ws.Range("E24").Value = '  +5.94%  '

$ws.Range("D25").Value = '''2.13'
$ws.Range("E25").Value = '  +11.98%  '

$ws.Range("E26").Value = '  +0.19%  '

$ws.Range("D27").Value = '''42.19'
$ws.Range("E27").Value = '  +3.52%  '

$ws.Range("D28").Value = '''24.24'
$ws.Range("E28").Value = '  +15.32%  '

$ws.Range("D29").Value = '''10.46'
$ws.Range("E29").Value = '  +7.42%  '

$ws.Range("D30").Value = '''2.27'
$ws.Range("E30").Value = '  +1.99%  '

$ws.Range("D31").Value = '''3.83'
$ws.Range("E31").Value = '  +4.29%  '

$ws.Range("D32").Value = '''6.02'
$ws.Range("E32").Value = '  +10.82%  '

$ws.Range("E33").Value = '  +6.41%  '

$ws.Range("D34").Value = '''0.0848'
$ws.Range("E34").Value = '  +9.40%  '

$ws.Range("D35").Value = '''2.19'
$ws.Range("E35").Value = '  +20.93%  '

$ws.Range("D36").Value = '''150.24'
$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("D37").Value = '''0.118'
$ws.Range("E37").Value = '  +3.72%  '

$ws.Range("D38").Value = '''0.122'
$ws.Range("E38").Value = '  +4.62%  '

$ws.Range("D39").Value = '''16.43'
$ws.Range("E39").Value = '  +7.64%  '

$ws.Range("D40").Value = '''4.28'
$ws.Range("E40").Value = '  +9.51%  '

$ws.Range("D41").Value = '''0.0329'
$ws.Range("E41").Value = '  +9.96%  '

$ws.Range("E42").Value = '  +11.00%  '

$ws.Range("D43").Value = '2.009.14'
$ws.Range("E43").Value = '  +8.03%  '

$ws.Range("D44").Value = '''0.999'
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").Value = '''94.19'
$ws.Range("E45").Value = '  +3.54%  '

$ws.Range("D46").Value = '''17.64'
$ws.Range("E46").Value = '  +36.80%  '

$ws.Range("D47").Value = '''1.89'
$ws.Range("E47").Value = '  +6.93%  '

$ws.Range("E48").Value = '  +9.57%  '

$ws.Range("D49").Value = '''8.98'
$ws.Range("E49").Value = '  +11.91%  '

$ws.Range("D50").Value = '''107.27'
$ws.Range("E50").Value = '  +11.08%  '

$ws.Range("B51").Value = 'ordi'
$ws.Range("C51").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D51").Value = '''74.94'
$ws.Range("E51").Value = '  +6.51%  '
